# Add data for 2021-12-27
# Updates the "through December 18" snapshot to "through December 19",
# and bumps several neighborhood/month counts that changed as new
# carjacking records were added for the additional day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date.
$ws.Name = "Through 2021-12-19"

# Update the running-month column header text (shared string used by B1).
$ws.Range("B1").Value = "December 2021 (through December 19)"

# Cell value updates (neighborhood row x month column).
$ws.Range("B3").Value = 7     # Englewood / December 2021 (through Dec 19)
$ws.Range("AX4").Value = 5    # North Lawndale / December 2017
$ws.Range("AX7").Value = 6    # Austin / December 2017
$ws.Range("BJ7").Value = 6    # Austin / December 2016
$ws.Range("BV7").Value = 3    # Austin / December 2015
$ws.Range("BJ11").Value = 4   # Humboldt Park / December 2016
$ws.Range("N13").Value = 2    # Roseland / December 2020
$ws.Range("N15").Value = 5    # Washington Heights / December 2020
$ws.Range("N23").Value = 2    # Little Village / December 2020
$ws.Range("AX24").Value = 5   # South Shore / December 2017
$ws.Range("N29").Value = 2    # Avalon Park / December 2020
$ws.Range("B33").Value = 3    # Near South Side / December 2021 (through Dec 19)
$ws.Range("N34").Value = 2    # Woodlawn / December 2020
$ws.Range("Z38").Value = 2    # Auburn Gresham / December 2019
$ws.Range("B53").Value = 1    # Hyde Park / December 2021 (through Dec 19) - new cell
$ws.Range("BJ73").Value = 2   # Fuller Park / December 2016
$ws.Range("N97").Value = 2    # Streeterville / December 2020
